$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("BF2:BF31")
$range.NumberFormat = "@"

for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 58).Value = "2014-06-22"
}
